$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2 (shifts old row2 data down to row3)
$ws.Rows.Item(2).Insert()

# New row 2: partial entry (id=2, nome=a, cognome=a, email=a@a.it, file=Screenshot 2024-08-30 alle 15.56.50.png)
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = "a"
$ws.Cells.Item(2, 3).Value = "a"
$ws.Cells.Item(2, 4).Value = "a@a.it"
$ws.Cells.Item(2, 5).Value = "Screenshot 2024-08-30 alle 15.56.50.png"

# Update row 3 (previously row 2) values: nome/cognome -> a, email -> a@a.it, file -> new screenshot name
$ws.Cells.Item(3, 2).Value = "a"
$ws.Cells.Item(3, 3).Value = "a"
$ws.Cells.Item(3, 4).Value = "a@a.it"
$ws.Cells.Item(3, 5).Value = "Screenshot 2024-10-21 alle 20.27.06.png"
